$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above row 11 ("Description") to make room for "Jurisdiction",
# pushing "Description" ... "Count" down from rows 11-21 to rows 12-22.
$ws.Rows.Item(11).Insert()

# The inserted row picks up a blank default format; copy the (unchanged)
# formatting of the row directly below it back onto the new row so it matches
# every other data row in the table.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# New "Jurisdiction" row (value left blank).
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update the "Date" value (row 8, column B).
$ws.Cells.Item(8, 2).Value = "2024-10-02T15:04:17+00:00"

# Update the "Contact" value (row 10, column B).
$ws.Cells.Item(10, 2).Value = "Ferlab.bio (http://example.org/example-publisher)"
